$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the first 16 data rows (rows 2-17), which hold the oldest
# quarters of data. This shifts all subsequent rows up by 16 and
# shrinks the used range from A1:B164 to A1:B148.
$ws.Range("A2:B17").EntireRow.Delete()
